# Remove the CF row for "Occupation, water courses, artificial".
# This string does not exist in the ecoinvent nomenclature (it should map to
# "Occupation, river, artificial" and "Occupation, lakes, artificial" instead).
# The corresponding data row in the sheet is row 32 - delete it entirely so
# every row below shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "Occupation, water courses, artificial"

$found = $ws.Columns.Item(1).Find($target)
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

$excel.ActiveWindow.ScrollRow = 17
$ws.Range("A35").Select()
